$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.278.71'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '2.587.32'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'572.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.77%  '
$ws.Range('D6').Value = "'143.69"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.597.16'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('E12').Value = '  +10.21%  '
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').Value = '3.042.38'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '59.305.87'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = "'22.55"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.19%  '
$ws.Range('E17').Value = '  +4.34%  '
$ws.Range('D18').Value = '2.590.12'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').Value = "'338.58"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'64.40"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.60%  '
$ws.Range('D25').Value = "'0.454"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.33%  '
$ws.Range('D26').Value = "'0.998"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('D28').Value = "'7.28"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.18%  '
$ws.Range('E29').Value = '  +3.96%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').Value = "'159.02"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').Value = "'0.879"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = "'0.883"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.88%  '
$ws.Range('D39').Value = "'37.09"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('D41').Value = "'294.43"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.83%  '
$ws.Range('E42').Value = '  +2.16%  '
$ws.Range('D43').Value = "'0.998"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = "'0.0977"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('D45').Value = "'0.597"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').Value = "'19.26"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.36%  '
$ws.Range('D48').Value = "'10.63"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = "'124.35"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.97%  '
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('D51').Value = '1.946.31'
